$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B and D act as thin "spacer" columns formatted with a bold,
# shaded, centered style (already applied to existing cells above).
# Extend that spacer formatting/content down to match how far the
# adjacent content columns (A and C) now run.

# Grab a template cell for each spacer column that already has the
# desired style applied, then clone its format into the new cells.
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B10").Value = "x"

$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B17").Value = "x"

$ws.Range("D13").Copy()
$ws.Range("D14:D26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D14").Value = "x"
$ws.Range("D15").Value = "x"
$ws.Range("D16").Value = "x"
$ws.Range("D17").Value = "x"
$ws.Range("D18").Value = "x"
$ws.Range("D19").Value = "x"
$ws.Range("D20").Value = "x"
$ws.Range("D21").Value = "x"
$ws.Range("D22").Value = "x"
$ws.Range("D23").Value = "x"
$ws.Range("D24").Value = "x"
$ws.Range("D25").Value = "x"
$ws.Range("D26").Value = "x"

$ws.Range("D13").Copy()
$ws.Range("D27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D27").Value = "xx"

$excel.CutCopyMode = 0

# Leave the active selection on the last edited cell.
$ws.Range("D27").Select()
